$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$ws.Range("K2").Value = "2025-11-21T07:01:52.039165+00:00"
$ws.Range("K3").Value = "2025-11-21T07:01:52.039197+00:00"
$ws.Range("K4").Value = "2025-11-21T07:01:54.321165+00:00"
$ws.Range("K5").Value = "2025-11-21T07:01:54.321192+00:00"
$ws.Range("K6").Value = "2025-11-21T07:01:56.589072+00:00"
$ws.Range("K7").Value = "2025-11-21T07:01:59.248078+00:00"
$ws.Range("K8").Value = "2025-11-21T07:02:01.407983+00:00"
$ws.Range("K9").Value = "2025-11-21T07:02:01.408011+00:00"
$ws.Range("K10").Value = "2025-11-21T07:02:01.408031+00:00"
$ws.Range("K11").Value = "2025-11-21T07:02:03.653028+00:00"
$ws.Range("K12").Value = "2025-11-21T07:02:05.895929+00:00"
$ws.Range("K13").Value = "2025-11-21T07:02:08.698251+00:00"
$ws.Range("K14").Value = "2025-11-21T07:02:11.344156+00:00"
$ws.Range("K15").Value = "2025-11-21T07:02:11.344185+00:00"
$ws.Range("K16").Value = "2025-11-21T07:02:11.344201+00:00"
$ws.Range("K17").Value = "2025-11-21T07:02:13.562906+00:00"
$ws.Range("K18").Value = "2025-11-21T07:02:18.824051+00:00"
$ws.Range("K19").Value = "2025-11-21T07:02:18.824081+00:00"
$ws.Range("K20").Value = "2025-11-21T07:02:21.078373+00:00"
$ws.Range("K21").Value = "2025-11-21T07:02:21.078404+00:00"
$ws.Range("K22").Value = "2025-11-21T07:02:21.078420+00:00"
$ws.Range("K23").Value = "2025-11-21T07:02:23.725915+00:00"
$ws.Range("K24").Value = "2025-11-21T07:02:23.725943+00:00"
$ws.Range("K25").Value = "2025-11-21T07:02:26.316174+00:00"
$ws.Range("K26").Value = "2025-11-21T07:02:26.316202+00:00"
$ws.Range("K27").Value = "2025-11-21T07:02:26.316220+00:00"
$ws.Range("K28").Value = "2025-11-21T07:02:28.975965+00:00"
$ws.Range("K29").Value = "2025-11-21T07:02:28.975994+00:00"
$ws.Range("K30").Value = "2025-11-21T07:02:31.550064+00:00"
$ws.Range("K31").Value = "2025-11-21T07:02:31.550093+00:00"
$ws.Range("K32").Value = "2025-11-21T07:02:31.550112+00:00"
$ws.Range("K33").Value = "2025-11-21T07:02:31.550127+00:00"
$ws.Range("K34").Value = "2025-11-21T07:02:31.550142+00:00"
$ws.Range("K35").Value = "2025-11-21T07:02:33.835510+00:00"
$ws.Range("K36").Value = "2025-11-21T07:02:33.835538+00:00"
$ws.Range("K37").Value = "2025-11-21T07:02:39.145985+00:00"
$ws.Range("K38").Value = "2025-11-21T07:02:39.146015+00:00"
$ws.Range("K39").Value = "2025-11-21T07:02:39.146044+00:00"
$ws.Range("K40").Value = "2025-11-21T07:02:41.375942+00:00"
